$p = $ppt.ActivePresentation

# 1) Nudge the "Picture 4" image on slide 10 ("What is OpenPose?") to its
#    new position. OOXML offsets (EMU) 6614610/1480494 -> 6619771/1501828,
#    converted to points (1 pt = 12700 EMU) for the COM Left/Top properties.
$s10 = $p.Slides.Item(10)
$pic = $s10.Shapes.Item(9)
$pic.Left = 521.2418110236221
$pic.Top = 118.25417322834646

# 2) Remove the trailing "Let's see it in action:" slide (the last slide in
#    the deck) - it was dropped from the presentation.
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Delete()
